$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.439551472663879
$ws.Range("B1").Value = 1.338711977005005
$ws.Range("C1").Value = 5.449761390686035
$ws.Range("D1").Value = 1.481338500976562
$ws.Range("E1").Value = 0.9833788275718689
